$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the three email addresses (add a digit before the @) while
# keeping their existing hyperlink/formatting untouched.
$ws.Range("C2").Value = "janfaizi1@gmail.com"
$ws.Range("C3").Value = "alijan2@tek.com"
$ws.Range("C4").Value = "anisa2@gmail.com"

# Move the active cell selection from F7 to D7.
$ws.Range("D7").Select()
